# Update export to display multiple geologic ages as a comma delimited list.
#
# The worksheet has a "Geologic Age Code" column (W) and a duplicate/auxiliary
# "Geologic Age Code" column (AR) that was used to hold a second age code for
# a sample. This script merges the AR value into the W column as a
# comma-delimited list ("<primary>, <secondary>") for every data row, then
# clears out the now-redundant AR column (including its header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 6

for ($r = 2; $r -le $lastRow; $r++) {
    $primary = $ws.Range("W$r").Value2
    $secondary = $ws.Range("AR$r").Value2

    if ($secondary -ne $null -and $secondary -ne "") {
        $ws.Range("W$r").Value = $primary.ToString() + ", " + $secondary.ToString()
    }

    $ws.Range("AR$r").ClearContents()
}

# The AR column no longer has a header either.
$ws.Range("AR1").ClearContents()
